# COP19_man_usa.xlsx update
# Re-codes several rows of the qualitative-coding sheet:
#  - Row 13 (policy target / emissions coding) is no longer coded as relevant -> clear B:H, set B="no"
#  - Row 14 (climate finance coding) is no longer coded as relevant -> clear B:H, set B="no"
#  - Row 16 "30 word explanation" (H) rewritten
#  - Row 17 Unit/Scale values rewritten, Principle + explanation rewritten
#  - Row 20 Principle + explanation rewritten
#  - Row 24 (complementary initiatives paragraph) is now coded as relevant -> fill B:H
#  - Row 27 explanation (H) rewritten

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: no longer relevant - clear all coding columns, mark B as "no"
$ws.Range("C13:H13").Clear()
$ws.Range("B13").Value = "no"

# Row 14: no longer relevant - clear all coding columns, mark B as "no"
$ws.Range("C14:H14").Clear()
$ws.Range("B14").Value = "no"

# The following assignments are ordered to match the order the new coding
# values were first introduced (keeps the shared-string table ordering
# consistent with the source edit).

# Row 16: update the explanation text
$ws.Range("H16").Value = "Prescribes the moral judgement that all countries should participate in the new agreement. "

# Row 17: explanation
$ws.Range("H17").Value = "Moral judgement on the need of accounting for self-differentiation and accounting for CBDR. "

# Row 20: principle
$ws.Range("G20").Value = "general normative statement"

# Row 17: principle
$ws.Range("G17").Value = "libertarian "

# Row 20: explanation
$ws.Range("H20").Value = "Value judgement on the need to have evolving categories to determine responsiblities. No emphasis on specific distribution."

# Row 27: explanation
$ws.Range("H27").Value = "Moral judgement of the need for cooperation in light of mutual responsibility. "

# Row 24: now coded as relevant - populate coding columns
$ws.Range("B24").Value = "yes"
$ws.Range("C24").Value = "new UNFCCC policy, other(complementary initiatives), mitigation"
$ws.Range("D24").Value = "measures"
$ws.Range("E24").Value = "global"
$ws.Range("F24").Value = "n.a."
$ws.Range("G24").Value = "general normative statement"
$ws.Range("H24").Value = "Statement on the need for cooperation to take on action to include more initiatives. No specific distribution. "

# Row 17: topic / unit (new unique topic string introduced last)
$ws.Range("C17").Value = "new UNFCCC policy, CBDR, self-differentiation, mitigation"
$ws.Range("D17").Value = "n.a."

# Update the selection/view state to match the final edit location
$ws.Range("C27").Select()
